# Update the title font on slide 10 ("Mitigation Strategies") so the
# run uses "Gill Sans MT" for Latin, East Asian and Complex Script text
# instead of the MS PGothic East Asian override.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$tr.Font.Name = "Gill Sans MT"
$tr.Font.NameFarEast = "Gill Sans MT"
$tr.Font.NameComplexScript = "Gill Sans MT"
